$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-01 Saturday" "2025-11-02 Sunday"

Replace-Text "66÷9=7, 3" "26÷7=3, 5"
Replace-Text "30÷6=5, 0" "68÷3=22, 2"
Replace-Text "99÷7=14, 1" "68÷6=11, 2"
Replace-Text "82÷7=11, 5" "56÷4=14, 0"
Replace-Text "85÷5=17, 0" "32÷8=4, 0"

Replace-Text "87÷5=17, 2" "67÷3=22, 1"
Replace-Text "30÷2=15, 0" "97÷4=24, 1"
Replace-Text "42÷2=21, 0" "55÷3=18, 1"
Replace-Text "66÷3=22, 0" "14÷7=2, 0"
Replace-Text "71÷3=23, 2" "89÷7=12, 5"

Replace-Text "52÷6=8, 4" "34÷8=4, 2"
Replace-Text "53÷6=8, 5" "47÷5=9, 2"
Replace-Text "52÷9=5, 7" "75÷7=10, 5"
Replace-Text "37÷8=4, 5" "46÷7=6, 4"
Replace-Text "19÷9=2, 1" "89÷6=14, 5"

Replace-Text "43÷7=6, 1" "65÷6=10, 5"
Replace-Text "12÷7=1, 5" "26÷7=3, 5"
Replace-Text "86÷6=14, 2" "81÷7=11, 4"
Replace-Text "26÷4=6, 2" "79÷8=9, 7"
Replace-Text "66÷7=9, 3" "93÷2=46, 1"

Replace-Text "67÷2=33, 1" "20÷4=5, 0"
Replace-Text "60÷2=30, 0" "22÷8=2, 6"
Replace-Text "65÷7=9, 2" "15÷5=3, 0"
Replace-Text "81÷9=9, 0" "90÷3=30, 0"
Replace-Text "80÷6=13, 2" "48÷7=6, 6"
